$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default, unstyled data-row look) used to reset style
# index after forcing a numeric-looking value to be stored as text, so
# that cells keep their original (unstyled) formatting.
$refStyle = $ws.Range('B2').Style

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$TextValue
    )
    $c = $ws.Range($CellRef)
    # Force text storage so a numeric-looking string (e.g. "241.94")
    # isn't silently converted into a real number by Excel.
    $c.NumberFormat = '@'
    $c.Value = $TextValue
    # Restore the original (default) cell style so no stray formatting
    # is introduced by the temporary text number-format.
    $c.Style = $refStyle
}

$ws.Range('D2').Value = '40.892.41'
$ws.Range('E2').Value = '  -6.95%  '
$ws.Range('D3').Value = '2.195.24'
$ws.Range('E3').Value = '  -7.41%  '
$ws.Range('E4').Value = '  -0.27%  '
Set-TextValue 'D5' '241.94'
$ws.Range('E5').Value = '  +0.45%  '
Set-TextValue 'D6' '0.622'
$ws.Range('E6').Value = '  -8.02%  '
Set-TextValue 'D7' '69.08'
$ws.Range('E7').Value = '  -7.24%  '
$ws.Range('E8').Value = '  +0.00%  '
Set-TextValue 'D9' '0.542'
$ws.Range('E9').Value = '  -12.28%  '
Set-TextValue 'D10' '0.0948'
$ws.Range('E10').Value = '  -7.74%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D11' '36.31'
$ws.Range('E11').Value = '  -3.78%  '
$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D12' '57.61'
$ws.Range('E12').Value = '  -4.86%  '
$ws.Range('E13').Value = '  -4.64%  '
Set-TextValue 'D14' '6.59'
$ws.Range('E14').Value = '  -10.10%  '
$ws.Range('D15').Value = '2.515.49'
$ws.Range('E15').Value = '  -7.66%  '
Set-TextValue 'D16' '14.64'
$ws.Range('E16').Value = '  -11.10%  '
Set-TextValue 'D17' '0.831'
$ws.Range('E17').Value = '  -10.36%  '
$ws.Range('D18').Value = '2.195.53'
$ws.Range('E18').Value = '  -7.42%  '
$ws.Range('D19').Value = '40.758.91'
$ws.Range('E19').Value = '  -7.27%  '
$ws.Range('D20').Value = '0.0₃0939'
$ws.Range('E20').Value = '  -9.45%  '
Set-TextValue 'D21' '72.48'
$ws.Range('E21').Value = '  -7.32%  '
Set-TextValue 'D22' '6.04'
$ws.Range('E22').Value = '  -8.47%  '
Set-TextValue 'D23' '229.52'
$ws.Range('E23').Value = '  -9.82%  '
Set-TextValue 'D24' '2.02'
$ws.Range('E24').Value = '  +7.35%  '
$ws.Range('E25').Value = '  +0.00%  '
Set-TextValue 'D26' '3.59'
$ws.Range('E26').Value = '  -4.80%  '
Set-TextValue 'D27' '2.40'
$ws.Range('E27').Value = '  -4.40%  '
$ws.Range('E28').Value = '  -5.05%  '
Set-TextValue 'D29' '9.65'
$ws.Range('E29').Value = '  -8.78%  '
Set-TextValue 'D30' '168.95'
$ws.Range('E30').Value = '  -4.13%  '
Set-TextValue 'D31' '20.22'
$ws.Range('E31').Value = '  -9.98%  '
$ws.Range('E32').Value = '  -9.93%  '
Set-TextValue 'D33' '0.123'
$ws.Range('E33').Value = '  -8.36%  '
Set-TextValue 'D34' '0.0699'
$ws.Range('E34').Value = '  -7.89%  '
Set-TextValue 'D35' '5.12'
$ws.Range('E35').Value = '  -5.57%  '
Set-TextValue 'D36' '4.57'
$ws.Range('E36').Value = '  -10.55%  '
Set-TextValue 'D37' '3.83'
$ws.Range('E37').Value = '  -0.24%  '
Set-TextValue 'D38' '23.81'
$ws.Range('E38').Value = '  +14.75%  '
Set-TextValue 'D39' '2.26'
$ws.Range('E39').Value = '  -6.79%  '
Set-TextValue 'D40' '0.0270'
$ws.Range('E40').Value = '  -3.88%  '
Set-TextValue 'D41' '5.77'
$ws.Range('E41').Value = '  -13.33%  '
Set-TextValue 'D42' '62.46'
$ws.Range('E42').Value = '  -4.17%  '
Set-TextValue 'D43' '4.85'
$ws.Range('E43').Value = '  -10.13%  '
Set-TextValue 'D44' '8.57'
$ws.Range('E44').Value = '  -5.92%  '
Set-TextValue 'D45' '0.191'
$ws.Range('E45').Value = '  -6.25%  '
$ws.Range('E46').Value = '  +0.22%  '
Set-TextValue 'D47' '0.0980'
$ws.Range('E47').Value = '  -8.98%  '
Set-TextValue 'D48' '4.49'
$ws.Range('E48').Value = '  +1.69%  '
Set-TextValue 'D49' '10.41'
$ws.Range('E49').Value = '  +6.67%  '
Set-TextValue 'D50' '1.16'
$ws.Range('E50').Value = '  -7.09%  '
Set-TextValue 'D51' '1.09'
$ws.Range('E51').Value = '  -6.52%  '
